$d = $word.ActiveDocument

# Locate the "Test4-" and "Test5-" paragraphs by their exact (current) text,
# then append the new trailing text to each, placing the insertion just
# before the paragraph mark so it lands inside the paragraph itself.

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()

    if ($t -eq "Test4-") {
        $target = $d.Range($p.Range.Start, $p.Range.End - 1)
        $target.Collapse(0)
        $target.InsertAfter(" ")
        $target.Collapse(0)
        $target.InsertAfter("https://www.hackerrank.com/contests/sda-2021-2022-test4-nov30/challenges")
    }
    elseif ($t -eq "Test5-") {
        $target = $d.Range($p.Range.Start, $p.Range.End - 1)
        $target.Collapse(0)
        $target.InsertAfter(" ")
    }
}
